$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume figures (GitHub Actions refresh).
# D-column "Price" cells are plain text in the source sheet (thousand-dot
# formatting like "42.807.82" is not a valid Excel number). Several new values
# (e.g. "253.93") DO parse as a number, so Excel would silently coerce them on
# assignment. Force text format before writing, then restore the default style
# so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.807.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.257.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.03"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.22%  "
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("E9").Value = "  +13.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0957"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.105"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.593.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.889"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.266.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.782.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0978"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.62%  "
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.129"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0785"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.74%  "
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "29.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("E40").Value = "  +8.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.202"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("E51").Value = "  +1.36%  "
